$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Material Dimensions & Physical Properties'
$ws.Range("B2").Value = 'DECIMALS'
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 'Parameter for logic calculation'
$ws.Range("A3").Value = 'Material Dimensions & Physical Properties'
$ws.Range("B3").Value = 'GLASS_DENSITY_KG_M3'
$ws.Range("C3").Value = 2500
$ws.Range("D3").Value = 'Density of glass (kg/m3)'
$ws.Range("A4").Value = 'Material Dimensions & Physical Properties'
$ws.Range("B4").Value = 'MASS_PER_M2_DOUBLE'
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 'Parameter for logic calculation'
$ws.Range("A5").Value = 'Material Dimensions & Physical Properties'
$ws.Range("B5").Value = 'MASS_PER_M2_SINGLE'
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 'Parameter for logic calculation'
$ws.Range("A6").Value = 'Material Dimensions & Physical Properties'
$ws.Range("B6").Value = 'MASS_PER_M2_TRIPLE'
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 'Parameter for logic calculation'
$ws.Range("A7").Value = 'Material Dimensions & Physical Properties'
$ws.Range("B7").Value = 'SEALANT_DENSITY_KG_M3'
$ws.Range("C7").Value = 1400
$ws.Range("D7").Value = 'Density of sealant (kg/m3)'
$ws.Range("A8").Value = 'Material Dimensions & Physical Properties'
$ws.Range("B8").Value = 'SPACER_MASS_PER_M_KG'
$ws.Range("C8").Value = 0.0648
$ws.Range("D8").Value = 'Linear mass of spacer (kg/m)'
$ws.Range("A9").Value = 'Material Emission Factors'
$ws.Range("B9").Value = 'EF_MAT_GLASS_100RC'
$ws.Range("C9").Value = 0.77
$ws.Range("D9").Value = 'Embodied Carbon: 100% Recycled Glass (kgCO2e/kg)'
$ws.Range("A10").Value = 'Material Emission Factors'
$ws.Range("B10").Value = 'EF_MAT_GLASS_VIRGIN'
$ws.Range("C10").Value = 1.29
$ws.Range("D10").Value = 'Embodied Carbon: Virgin Float Glass (kgCO2e/kg)'
$ws.Range("A11").Value = 'Material Emission Factors'
$ws.Range("B11").Value = 'EF_MAT_PVB'
$ws.Range("C11").Value = 4.683
$ws.Range("D11").Value = 'Embodied Carbon: PVB Interlayer (kgCO2e/kg)'
$ws.Range("A12").Value = 'Material Emission Factors'
$ws.Range("B12").Value = 'EF_MAT_SEALANT'
$ws.Range("C12").Value = 3.51
$ws.Range("D12").Value = 'Embodied Carbon: Generic Sealant (kgCO2e/kg)'
$ws.Range("A13").Value = 'Material Emission Factors'
$ws.Range("B13").Value = 'EF_MAT_SPACER_ALU'
$ws.Range("C13").Value = 0.57
$ws.Range("D13").Value = 'Embodied Carbon: Aluminium Spacer (kgCO2e/linearmetre)'
$ws.Range("A14").Value = 'Material Emission Factors'
$ws.Range("B14").Value = 'EF_MAT_SPACER_STEEL'
$ws.Range("C14").Value = 0.17
$ws.Range("D14").Value = 'Embodied Carbon: Aluminium Spacer (kgCO2e/linearmetre)'
$ws.Range("A15").Value = 'Material Emission Factors'
$ws.Range("B15").Value = 'EF_MAT_SPACER_SWISS'
$ws.Range("C15").Value = 0.14
$ws.Range("D15").Value = 'Embodied Carbon: Aluminium Spacer (kgCO2e/linearmetre)'
$ws.Range("A16").Value = 'Miscellaneous'
$ws.Range("B16").Value = 'EF_PROCESS_COATING'
$ws.Range("C16").Value = 0.27
$ws.Range("D16").Value = 'Parameter for logic calculation'
$ws.Range("A17").Value = 'Miscellaneous'
$ws.Range("B17").Value = 'EF_PROCESS_LAMINATING'
$ws.Range("C17").Value = 0.3
$ws.Range("D17").Value = 'Parameter for logic calculation'
$ws.Range("A18").Value = 'Miscellaneous'
$ws.Range("B18").Value = 'EF_PROCESS_TOUGHENING'
$ws.Range("C18").Value = 0.39
$ws.Range("D18").Value = 'Parameter for logic calculation'
$ws.Range("A19").Value = 'Miscellaneous'
$ws.Range("B19").Value = 'MAX_TRUCK_LOAD_KG'
$ws.Range("C19").Value = 24000
$ws.Range("D19").Value = 'Parameter for logic calculation'
$ws.Range("A20").Value = 'Process Emission Factors'
$ws.Range("B20").Value = 'BREAKING_KGCO2_PER_M2'
$ws.Range("C20").Value = 0.05
$ws.Range("D20").Value = 'Parameter for logic calculation'
$ws.Range("A21").Value = 'Process Emission Factors'
$ws.Range("B21").Value = 'DISASSEMBLY_KGCO2_PER_M2'
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 'Emissions for disassembly process per m2'
$ws.Range("A22").Value = 'Process Emission Factors'
$ws.Range("B22").Value = 'E_SITE_KGCO2_PER_M2'
$ws.Range("C22").Value = 0.15
$ws.Range("D22").Value = 'Site energy emissions per m2'
$ws.Range("A23").Value = 'Process Emission Factors'
$ws.Range("B23").Value = 'INSTALL_SYSTEM_KGCO2_PER_M2'
$ws.Range("C23").Value = 0.15
$ws.Range("D23").Value = 'Parameter for logic calculation'
$ws.Range("A24").Value = 'Process Emission Factors'
$ws.Range("B24").Value = 'PROCESS_ENERGY_ASSEMBLY_KGCO2_PER_M2'
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 'Energy emission for IGU assembly per m2'
$ws.Range("A25").Value = 'Process Emission Factors'
$ws.Range("B25").Value = 'RECONDITION_KGCO2_PER_M2'
$ws.Range("C25").Value = 0.5
$ws.Range("D25").Value = 'Parameter for logic calculation - used for reconditioning components'
$ws.Range("A26").Value = 'Process Emission Factors'
$ws.Range("B26").Value = 'REMANUFACTURING_KGCO2_PER_M2'
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 'Emissions for remanufacturing process per m2'
$ws.Range("A27").Value = 'Process Emission Factors'
$ws.Range("B27").Value = 'REPAIR_KGCO2_PER_M2'
$ws.Range("C27").Formula = "=0.596*0.03"
$ws.Range("C27").NumberFormat = "0.00"
$ws.Range("D27").Value = 'Parameter for logic calculation (kgCO2/cavity Argon fill)'
$ws.Range("A28").Value = 'Process Emission Factors'
$ws.Range("B28").Value = 'REPURPOSE_HEAVY_KGCO2_PER_M2'
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 'Parameter for logic calculation'
$ws.Range("A29").Value = 'Process Emission Factors'
$ws.Range("B29").Value = 'REPURPOSE_LIGHT_KGCO2_PER_M2'
$ws.Range("C29").Value = 0.5
$ws.Range("D29").Value = 'Parameter for logic calculation'
$ws.Range("A30").Value = 'Process Emission Factors'
$ws.Range("B30").Value = 'REPURPOSE_MEDIUM_KGCO2_PER_M2'
$ws.Range("C30").Value = 0.75
$ws.Range("D30").Value = 'Parameter for logic calculation'
$ws.Range("A31").Value = 'Process Yields & Losses'
$ws.Range("B31").Value = 'BREAKAGE_RATE_GLOBAL'
$ws.Range("C31").Value = 0.01
$ws.Range("D31").Value = 'Global breakage rate estimate (0.0-1.0)'
$ws.Range("A32").Value = 'Process Yields & Losses'
$ws.Range("B32").Value = 'HUMIDITY_FAILURE_RATE'
$ws.Range("C32").Value = 0.02
$ws.Range("D32").Value = 'Rate of units failing humidity tests (0.0-1.0)'
$ws.Range("A33").Value = 'Process Yields & Losses'
$ws.Range("B33").Value = 'REMANUFACTURING_YIELD'
$ws.Range("C33").Value = 0.95
$ws.Range("D33").Value = 'Yield of remanufacturing acceptable components (0.0-1.0)'
$ws.Range("A34").Value = 'Process Yields & Losses'
$ws.Range("B34").Value = 'SPLIT_YIELD'
$ws.Range("C34").Value = 0.9
$ws.Range("D34").Value = 'Success rate of splitting panes (0.0-1.0)'
$ws.Range("A35").Value = 'Process Yields & Losses'
$ws.Range("B35").Value = 'YIELD_SYSTEM_REUSE'
$ws.Range("C35").Value = 0.2
$ws.Range("D35").Value = 'Yield loss associated with system reuse (0.0-1.0)'
$ws.Range("A36").Value = 'Process Yields & Losses'
$ws.Range("B36").Value = 'YIELD_REPAIR'
$ws.Range("C36").Value = 0.1
$ws.Range("D36").Value = 'Yield loss during repair process for system repair (0.0-1.0)'
$ws.Range("A37").Value = 'Process Yields & Losses'
$ws.Range("B37").Value = 'YIELD_DISASSEMBLY_REUSE'
$ws.Range("C37").Value = 0.2
$ws.Range("D37").Value = 'Yield loss during disassembly for component reuse (0.0-1.0)'
$ws.Range("A38").Value = 'Process Yields & Losses'
$ws.Range("B38").Value = 'YIELD_DISASSEMBLY_REMANUFACTURE'
$ws.Range("C38").Value = 0.1
$ws.Range("D38").Value = 'Yield loss during disassembly for remanufacture (0.0-1.0)'
$ws.Range("A39").Value = 'Process Yields & Losses'
$ws.Range("B39").Value = 'YIELD_DISASSEMBLY_REPURPOSE'
$ws.Range("C39").Value = 0.1
$ws.Range("D39").Value = 'Yield loss during disassembly for repurpose (0.0-1.0)'
$ws.Range("A40").Value = 'Recycling & Credits'
$ws.Range("B40").Value = 'SHARE_CULLET_FLOAT'
$ws.Range("C40").Value = 0.8
$ws.Range("D40").Value = 'Parameter for logic calculation'
$ws.Range("A41").Value = 'Recycling & Credits'
$ws.Range("B41").Value = 'SHARE_CULLET_OPEN_LOOP_CONT'
$ws.Range("C41").Value = 0.4
$ws.Range("D41").Value = 'Parameter for logic calculation'
$ws.Range("A42").Value = 'Recycling & Credits'
$ws.Range("B42").Value = 'SHARE_CULLET_OPEN_LOOP_GW'
$ws.Range("C42").Value = 0.4
$ws.Range("D42").Value = 'Parameter for logic calculation'
$ws.Range("A43").Value = 'Recycling & Credits'
$ws.Range("B43").Value = 'FLOAT_GLASS_REPROCESSING'
$ws.Range("C43").Value = 0.52
$ws.Range("D43").Value = 'Embodied Carbon: Reprocessing Flat Glass (kgCO2e/kg)'
$ws.Range("A44").Value = 'Stillage & Logistics'
$ws.Range("B44").Value = 'IGUS_PER_STILLAGE'
$ws.Range("C44").Value = 10
$ws.Range("D44").Value = 'Parameter for logic calculation'
$ws.Range("A45").Value = 'Stillage & Logistics'
$ws.Range("B45").Value = 'INCLUDE_STILLAGE_EMBODIED'
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 'Parameter for logic calculation'
$ws.Range("A46").Value = 'Stillage & Logistics'
$ws.Range("B46").Value = 'STILLAGE_LIFETIME_CYCLES'
$ws.Range("C46").Value = 1000
$ws.Range("D46").Value = 'Parameter for logic calculation'
$ws.Range("A47").Value = 'Stillage & Logistics'
$ws.Range("B47").Value = 'STILLAGE_MANUFACTURE_KGCO2'
$ws.Range("C47").Value = 150
$ws.Range("D47").Value = 'Parameter for logic calculation'
$ws.Range("A48").Value = 'Stillage & Logistics'
$ws.Range("B48").Value = 'STILLAGE_MASS_EMPTY_KG'
$ws.Range("C48").Value = 80
$ws.Range("D48").Value = 'Parameter for logic calculation'
$ws.Range("A49").Value = 'System Configuration'
$ws.Range("B49").Value = 'GEOCODER_USER_AGENT'
$ws.Range("C49").Value = 'igu-reuse-tool/0.1 (CHANGE_THIS_TO_YOUR_EMAIL@DOMAIN)'
$ws.Range("D49").Value = 'Parameter for logic calculation'
$ws.Range("A50").Value = 'Transport Settings'
$ws.Range("B50").Value = 'BACKHAUL_FACTOR'
$ws.Range("C50").Value = 1.6
$ws.Range("D50").Value = 'Backhaul adjustment factor (>1.0)'
$ws.Range("A51").Value = 'Transport Settings'
$ws.Range("B51").Value = 'DISTANCE_FALLBACK_A_KM'
$ws.Range("C51").Value = 50
$ws.Range("D51").Value = 'Default distance Origin->Processor (km)'
$ws.Range("A52").Value = 'Transport Settings'
$ws.Range("B52").Value = 'DISTANCE_FALLBACK_B_KM'
$ws.Range("C52").Value = 50
$ws.Range("D52").Value = 'Default distance Processor->Reuse (km)'
$ws.Range("A53").Value = 'Transport Settings'
$ws.Range("B53").Value = 'EMISSIONFACTOR_FERRY'
$ws.Range("C53").Value = 0.015
$ws.Range("D53").Value = 'Ferry emission factor (kgCO2e/tkm)'
$ws.Range("A54").Value = 'Transport Settings'
$ws.Range("B54").Value = 'EMISSIONFACTOR_TRUCK'
$ws.Range("C54").Value = 0.062
$ws.Range("D54").Value = 'Truck emission factor (kgCO2e/tkm)'
$ws.Range("A55").Value = 'Transport Settings'
$ws.Range("B55").Value = 'FERRY_CAPACITY_T'
$ws.Range("C55").Value = 500
$ws.Range("D55").Value = 'Parameter for logic calculation'
$ws.Range("A56").Value = 'Transport Settings'
$ws.Range("B56").Value = 'ROUTE_A_MODE'
$ws.Range("C56").Value = 'HGV lorry'
$ws.Range("D56").Value = 'Parameter for logic calculation'
$ws.Range("A57").Value = 'Transport Settings'
$ws.Range("B57").Value = 'ROUTE_B_MODE'
$ws.Range("C57").Value = 'HGV lorry'
$ws.Range("D57").Value = 'Parameter for logic calculation'
$ws.Range("A58").Value = 'Transport Settings'
$ws.Range("B58").Value = 'TRUCK_CAPACITY_T'
$ws.Range("C58").Value = 24
$ws.Range("D58").Value = 'Parameter for logic calculation'
$ws.Range("A59").Value = 'General Settings'
$ws.Range("B59").Value = 'Default IGU Service Lifetime (years)'
$ws.Range("C59").Value = 25
$ws.Range("D59").Value = 'Default age assumption for IGUs at end-of-life'

# Sheet view changes
$ws.Range("G15").Select()
$ws.Application.ActiveWindow.Zoom = 106

# Column A width (approximate match to target stored width 33.6640625)
$ws.Columns.Item(1).ColumnWidth = 32.83
